# Add team record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1) - set the new header text, then copy the existing
# header cell's formatting (bold, centered, bordered) onto the new cells
# so they match the rest of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-51: every team played 81 home games with a record of
# 81 wins, 81 losses, and 0 ties.
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 81  # AD
    $ws.Cells.Item($r, 31).Value = 81  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
